$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- Refresh "time_taken" (column F) timestamps on the data sheet ---
$ws1.Cells.Item(2, 6).Value = "2021-10-05 14:21:53.397402"
$ws1.Cells.Item(3, 6).Value = "2021-10-05 14:21:53.397410"
$ws1.Cells.Item(4, 6).Value = "2021-10-05 14:21:53.397413"
$ws1.Cells.Item(5, 6).Value = "2021-10-05 14:21:53.397416"
$ws1.Cells.Item(6, 6).Value = "2021-10-05 14:21:53.397419"
$ws1.Cells.Item(7, 6).Value = "2021-10-05 14:21:53.397422"
$ws1.Cells.Item(8, 6).Value = "2021-10-05 14:21:53.397424"
$ws1.Cells.Item(9, 6).Value = "2021-10-05 14:21:53.397427"
$ws1.Cells.Item(10, 6).Value = "2021-10-05 14:21:53.397429"
$ws1.Cells.Item(11, 6).Value = "2021-10-05 14:21:53.397432"
$ws1.Cells.Item(12, 6).Value = "2021-10-05 14:21:53.397434"
$ws1.Cells.Item(13, 6).Value = "2021-10-05 14:21:53.397437"
$ws1.Cells.Item(14, 6).Value = "2021-10-05 14:21:53.397439"
$ws1.Cells.Item(15, 6).Value = "2021-10-05 14:21:53.397442"
$ws1.Cells.Item(16, 6).Value = "2021-10-05 14:21:53.397444"
$ws1.Cells.Item(17, 6).Value = "2021-10-05 14:21:53.397447"
$ws1.Cells.Item(18, 6).Value = "2021-10-05 14:21:53.397449"
$ws1.Cells.Item(19, 6).Value = "2021-10-05 14:21:53.397452"
$ws1.Cells.Item(20, 6).Value = "2021-10-05 14:21:53.397454"
$ws1.Cells.Item(21, 6).Value = "2021-10-05 14:21:53.397457"
$ws1.Cells.Item(22, 6).Value = "2021-10-05 14:21:53.397460"
$ws1.Cells.Item(23, 6).Value = "2021-10-05 14:21:53.397462"
$ws1.Cells.Item(24, 6).Value = "2021-10-05 14:21:53.397465"
$ws1.Cells.Item(25, 6).Value = "2021-10-05 14:21:53.397467"
$ws1.Cells.Item(26, 6).Value = "2021-10-05 14:21:53.397470"
$ws1.Cells.Item(27, 6).Value = "2021-10-05 14:21:53.397473"
$ws1.Cells.Item(28, 6).Value = "2021-10-05 14:21:53.397475"
$ws1.Cells.Item(29, 6).Value = "2021-10-05 14:21:53.397478"
$ws1.Cells.Item(30, 6).Value = "2021-10-05 14:21:53.397480"
$ws1.Cells.Item(31, 6).Value = "2021-10-05 14:21:53.397482"
$ws1.Cells.Item(32, 6).Value = "2021-10-05 14:21:53.397485"
$ws1.Cells.Item(33, 6).Value = "2021-10-05 14:21:53.397487"
$ws1.Cells.Item(34, 6).Value = "2021-10-05 14:21:53.397490"
$ws1.Cells.Item(35, 6).Value = "2021-10-05 14:21:53.397493"
$ws1.Cells.Item(36, 6).Value = "2021-10-05 14:21:53.397495"
$ws1.Cells.Item(37, 6).Value = "2021-10-05 14:21:53.397498"
$ws1.Cells.Item(38, 6).Value = "2021-10-05 14:21:53.397500"
$ws1.Cells.Item(39, 6).Value = "2021-10-05 14:21:53.397503"
$ws1.Cells.Item(40, 6).Value = "2021-10-05 14:21:53.397505"
$ws1.Cells.Item(41, 6).Value = "2021-10-05 14:21:53.397508"
$ws1.Cells.Item(42, 6).Value = "2021-10-05 14:21:53.397511"
$ws1.Cells.Item(43, 6).Value = "2021-10-05 14:21:53.397513"
$ws1.Cells.Item(44, 6).Value = "2021-10-05 14:21:53.397516"
$ws1.Cells.Item(45, 6).Value = "2021-10-05 14:21:53.397518"
$ws1.Cells.Item(46, 6).Value = "2021-10-05 14:21:53.397521"
$ws1.Cells.Item(47, 6).Value = "2021-10-05 14:21:53.397523"
$ws1.Cells.Item(48, 6).Value = "2021-10-05 14:21:53.397526"
$ws1.Cells.Item(49, 6).Value = "2021-10-05 14:21:53.397529"
$ws1.Cells.Item(50, 6).Value = "2021-10-05 14:21:53.397531"
$ws1.Cells.Item(51, 6).Value = "2021-10-05 14:21:53.397533"
$ws1.Cells.Item(52, 6).Value = "2021-10-05 14:21:53.397536"
$ws1.Cells.Item(53, 6).Value = "2021-10-05 14:21:53.397538"
$ws1.Cells.Item(54, 6).Value = "2021-10-05 14:21:53.397541"
$ws1.Cells.Item(55, 6).Value = "2021-10-05 14:21:53.397544"
$ws1.Cells.Item(56, 6).Value = "2021-10-05 14:21:53.397547"
$ws1.Cells.Item(57, 6).Value = "2021-10-05 14:21:53.397549"
$ws1.Cells.Item(58, 6).Value = "2021-10-05 14:21:53.397552"
$ws1.Cells.Item(59, 6).Value = "2021-10-05 14:21:53.397554"
$ws1.Cells.Item(60, 6).Value = "2021-10-05 14:21:53.397557"
$ws1.Cells.Item(61, 6).Value = "2021-10-05 14:21:53.397559"
$ws1.Cells.Item(62, 6).Value = "2021-10-05 14:21:53.397562"
$ws1.Cells.Item(63, 6).Value = "2021-10-05 14:21:53.397564"
$ws1.Cells.Item(64, 6).Value = "2021-10-05 14:21:53.397567"
$ws1.Cells.Item(65, 6).Value = "2021-10-05 14:21:53.397570"
$ws1.Cells.Item(66, 6).Value = "2021-10-05 14:21:53.397573"

# --- Add the new "metadata" worksheet, positioned right after "data" ---
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "metadata"

# Header row
$new.Cells.Item(1, 2).Value = "data_name"
$new.Cells.Item(1, 3).Value = "data_id"
$new.Cells.Item(1, 4).Value = "data_version"
$new.Cells.Item(1, 5).Value = "data_version_created"
$new.Cells.Item(1, 6).Value = "panel_query_time"
$new.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$new.Cells.Item(2, 1).Value = 0
$new.Cells.Item(2, 2).Value = "Ocular coloboma"
$new.Cells.Item(2, 3).Value = 294
$new.Cells.Item(2, 4).NumberFormat = "@"
$new.Cells.Item(2, 4).Value = "1.44"
$new.Cells.Item(2, 4).ClearFormats()
$new.Cells.Item(2, 5).Value = "2021-09-09T10:36:56.195531Z"
$new.Cells.Item(2, 6).Value = "2021-10-05 14:21:53.394136"
$new.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/294/?format=json"

# Reuse the existing bold/centered/bordered header style from the data sheet
# (columns B:F) and the index-column style from A2, so no new style entries
# are introduced; then stamp the same style onto the extra G1 header cell.
$ws1.Range("B1:F1").Copy()
$new.Range("B1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$new.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)
